# chore: adapt column header formatting to respective input file names
#
#  - header cells that end in "_old" become "_FV2210" (the older input file's
#    format version), and header cells that end in "_new" become "_FV2304"
#    (the newer input file's format version)
#  - the sheet's data range is turned into a proper Excel Table ("Table1")
#  - the header row is frozen so it stays visible while scrolling

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header row (row 1) -------------------------------------
$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$firstCol = $usedRange.Column
$lastRow = $firstRow + $usedRange.Rows.Count - 1
$lastCol = $firstCol + $usedRange.Columns.Count - 1

for ($c = $firstCol; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item($firstRow, $c)
    $header = [string]$cell.Value2
    if ($header -like "*_old") {
        $cell.Value = ($header -replace "_old$", "_FV2210")
    } elseif ($header -like "*_new") {
        $cell.Value = ($header -replace "_new$", "_FV2304")
    }
}

# --- 2) Wrap the full data range (incl. the renamed header) in a Table ----
$dataRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# --- 3) Freeze the header row ----------------------------------------------
$ws.Activate()
$ws.Cells.Item($firstRow + 1, $firstCol).Select()
$excel.ActiveWindow.FreezePanes = $true
